$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 2) describing MCH144-1
$ws.Range("A2").Value = "MCH144-1"
$ws.Range("C2").Value = "HUMAN RIGHTS WATCH AFRICA REPORTS 1984-1994"
$ws.Range("D2").Value = "1984-1994"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21N | GRAP COUNT NUMER: NONE"

# Match the body font/color used for this new row (Calibri 10pt, automatic/theme text color)
# (column B is intentionally left blank/untouched, matching the source data)
foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.ThemeColor = 1
}

# Re-select the new row and keep the header row frozen, as in the saved view
$ws.Range("A2:J2").Select()
$excel.ActiveWindow.FreezePanes = $true
